$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update product/ingredient name cells to their de-duplicated / renumbered text values.
# (Backend refactor: product equality is now based on name+unit, so duplicate shared strings
#  that only differed by trailing whitespace / trailing comma / casing collapse to one value.)
# The first six assignments below introduce the brand-new consolidated names; they are set
# first (in this specific order) so the workbook's shared-string table lists them in the same
# order as the target file.
$ws.Range('E100').Value = 'ser zółty'
$ws.Range('E80').Value = 'schab bez kości'
$ws.Range('E16').Value = 'bułka kajzerka'
$ws.Range('E84').Value = 'ziemniaki'
$ws.Range('E102').Value = 'keczup'
$ws.Range('E104').Value = 'ryż'
$ws.Range('B1').Value = 'śniadanie'
$ws.Range('E1').Value = 'mleko'
$ws.Range('E3').Value = 'maliny'
$ws.Range('E4').Value = 'jeżyny'
$ws.Range('E5').Value = 'X'
$ws.Range('B6').Value = 'śniadanie'
$ws.Range('E6').Value = 'jogurt owocowy'
$ws.Range('E10').Value = 'X'
$ws.Range('B11').Value = 'śniadanie'
$ws.Range('E11').Value = 'jogurt owocowy'
$ws.Range('E15').Value = 'X'
$ws.Range('B16').Value = 'śniadanie'
$ws.Range('I16').Value = 'Bułkę kajzerkę przekroić na pół; wyciągnąć miąższ ze środka; posmarować cienko margaryną (cienko!) i włożyć do tostera; na patelni przygotować 2 małe jajka sadzone, posolić, popieprzyć; na grzankach w zagłębieniu położyć po plasterku szynki z kurczaka, na nich położyć jajko sadzone; całość obficie posypać szczypiorkiem; podawać z pomidorkami koktajlowymi.'
$ws.Range('E17').Value = 'jajko'
$ws.Range('E20').Value = 'pomidor'
$ws.Range('E21').Value = 'X'
$ws.Range('B22').Value = 'śniadanie'
$ws.Range('I22').Value = 'Rano wymieszać: sok z czarnej porzeczki, mleko 2% tł, jogurt naturalny, płatki owsiane, miód; schłodzić w lodówce.; przed zjedzeniem dodać orzechy włoskie'
$ws.Range('E23').Value = 'mleko'
$ws.Range('E24').Value = 'jogurt naturalny'
$ws.Range('H24').Value = '1/2 opako.'
$ws.Range('E25').Value = 'płatki owsiane '
$ws.Range('H25').Value = '2 kopiaste łyzki'
$ws.Range('E26').Value = 'miód '
$ws.Range('H26').Value = 'niepełna łyzka'
$ws.Range('E27').Value = 'orzechy włoskie '
$ws.Range('H27').Value = '1 łyzka'
$ws.Range('E28').Value = 'X'
$ws.Range('B29').Value = 'śniadanie'
$ws.Range('D29').Value = 'Jajecznica na masełku, pumpernikiel'
$ws.Range('E29').Value = 'jajko'
$ws.Range('H29').Value = '1 szt.'
$ws.Range('I29').Value = 'na patelni teflonowej na maśle wbić wymieszane jajko z białkami; doprawić solą i pieprzem ; dodać szczypiorek; zjeść z chlebem i pomidorem pokrojonym w ćwiartki. '
$ws.Range('E30').Value = 'jajko'
$ws.Range('H30').Value = 'tylko białko'
$ws.Range('E31').Value = 'pumpernikiel '
$ws.Range('H31').Value = '2 kromki'
$ws.Range('E32').Value = 'pomidor'
$ws.Range('H32').Value = '1 duzy'
$ws.Range('E33').Value = 'X'
$ws.Range('B34').Value = 'lunch'
$ws.Range('D34').Value = 'Chleb żytni razowy Do przegryzania wafle ryżowe'
$ws.Range('E34').Value = 'chleb żytni razowy'
$ws.Range('H34').Value = '1 kromka'
$ws.Range('I34').Value = 'Do przegryzania: wafle ryżowe słodkie ulubiony rodzaj (np. w polewie jogurtowej, 2 szt).'
$ws.Range('E35').Value = 'masło'
$ws.Range('H35').Value = '1/2 łyzeczki'
$ws.Range('H36').Value = '1 plasterek'
$ws.Range('E37').Value = 'sałata'
$ws.Range('H37').Value = '1 liść'
$ws.Range('E38').Value = 'wafle ryżowe'
$ws.Range('H38').Value = '2 szt.'
$ws.Range('E39').Value = 'jogurt owocowy'
$ws.Range('H39').Value = '1 szt.'
$ws.Range('E40').Value = 'X'
$ws.Range('B41').Value = 'lunch'
$ws.Range('D41').Value = 'Banan, Wafle ryżowe'
$ws.Range('E41').Value = 'banan'
$ws.Range('H41').Value = '2 szt.'
$ws.Range('I41').Value = 'Zjesc ze smakiem'
$ws.Range('E42').Value = 'wafle ryżowe'
$ws.Range('H42').Value = '3 szy.'
$ws.Range('E43').Value = 'X'
$ws.Range('B44').Value = 'lunch'
$ws.Range('D44').Value = 'Kanapki z serkiem topionym'
$ws.Range('E44').Value = 'bułka grahamka'
$ws.Range('H44').Value = '1 szt.'
$ws.Range('I44').Value = 'Zjesc ze smakiem'
$ws.Range('H45').Value = '1 łyzeczka'
$ws.Range('E46').Value = 'zółty ser'
$ws.Range('H46').Value = '1 plaster'
$ws.Range('E47').Value = 'sałata'
$ws.Range('H47').Value = '1 lisc'
$ws.Range('E48').Value = 'ogórek kiszony'
$ws.Range('H48').Value = '1 szt.'
$ws.Range('E49').Value = 'X'
$ws.Range('B50').Value = 'lunch'
$ws.Range('D50').Value = 'szybki lunch 1'
$ws.Range('E50').Value = 'bułka'
$ws.Range('H50').Value = '1 szt.'
$ws.Range('I50').Value = 'Malyszowe danie'
$ws.Range('K50').Value = 'Adam Malysz'
$ws.Range('E51').Value = 'banan'
$ws.Range('H51').Value = '2 szt.'
$ws.Range('E52').Value = 'X'
$ws.Range('B53').Value = 'lunch'
$ws.Range('D53').Value = 'szybki lunch 2'
$ws.Range('E53').Value = 'nutella'
$ws.Range('H53').Value = '1/2 opakowania'
$ws.Range('I53').Value = 'Malysz wersja na kwasno'
$ws.Range('K53').Value = 'Adam Malysz'
$ws.Range('E54').Value = 'sledzik'
$ws.Range('H54').Value = 'wiadro'
$ws.Range('E55').Value = 'X'
$ws.Range('B56').Value = 'lunch'
$ws.Range('D56').Value = 'szybki lunch 3'
$ws.Range('E56').Value = 'kapusta kwaszona'
$ws.Range('H56').Value = 'wiaderko'
$ws.Range('I56').Value = 'Przeklaska Kwasniewskiego'
$ws.Range('K56').Value = 'Kwaśniewski'
$ws.Range('E57').Value = 'X'
$ws.Range('B58').Value = 'obiad'
$ws.Range('D58').Value = 'Gulasz z kurczakiem i fasolką szparagową'
$ws.Range('E58').Value = 'filet z piersi kurczaka'
$ws.Range('H58').Value = '1 szt.'
$ws.Range('I58').Value = 'Filet z piersi kurczaka pokroić na kawałki; przyprawić solą, pieprzem i podsmażyć na patelni grillowej; pokrojone pieczarki poddusić z solą na oleju rzepakowym; dodać fasolkę szparagową i startą drobno marchewkę; dodać piersi z kurczaka, przyprawić ulubionymi ziołami; dodać natkę pietruszki. Podawać z ugotowaną kaszą gryczaną z jogurtem naturalnym wymieszanym z pietruszką.'
$ws.Range('E59').Value = 'pieczarki'
$ws.Range('H59').Value = '8 szt.'
$ws.Range('E60').Value = 'olej rzepakowy'
$ws.Range('H60').Value = '1 łyzka'
$ws.Range('E61').Value = 'fasolka szparagowa'
$ws.Range('H61').Value = '1/2 opakow.'
$ws.Range('E62').Value = 'marchew'
$ws.Range('H62').Value = '2 szt.'
$ws.Range('E63').Value = 'kasza gryczana'
$ws.Range('H63').Value = '3/4 woreczka'
$ws.Range('E64').Value = 'jogurt naturalny'
$ws.Range('H64').Value = '1/2 opakow.'
$ws.Range('E65').Value = 'X'
$ws.Range('B66').Value = 'obiad'
$ws.Range('D66').Value = 'Hamburgery drobiowe, Surówka z białej kapusty'
$ws.Range('E66').Value = 'biała kapusta'
$ws.Range('H66').Value = 'kolka liści'
$ws.Range('I66').Value = '1# przygotować surówkę:; białą kapustę pokroić cienko;  zmorzyć solą, odstawić na kilka minut; następnie dodać oliwę z oliwek;  doprawić solą, pieprzem i oregano;  2# Bułkę do hamburgerów [można zastąpić zwykłą bułką pszenną 60 g, średnia szt]; posmarować z dwóch stron sosem czosnkowym i keczupem; na jednej połówce ułożyć pokrojone warzywa:;  paprykę czerwoną; żółtą  i sałatę;  przygotować mięso: zmielone mięso z piersi z kurczaka; wymieszać z pokrojoną cebulą, białkiem jaja i namoczoną w wodzie bułką; doprawić do smaku solą, pieprzem i np. papryką ostrą; smażyć na patelni grillowej do rumianego koloru.'
$ws.Range('E67').Value = 'oliwa z oliwek'
$ws.Range('H67').Value = '1 łyżeczka'
$ws.Range('E68').Value = 'bułka do hamburgerów'
$ws.Range('H68').Value = '1 szt.'
$ws.Range('E69').Value = 'majonez light'
$ws.Range('H69').Value = '1 łyżeczka'
$ws.Range('E70').Value = 'jogurt naturalny'
$ws.Range('E71').Value = 'czosnek'
$ws.Range('H71').Value = '1 zabek'
$ws.Range('E72').Value = 'ketchup'
$ws.Range('H72').Value = '1 łyżeczka'
$ws.Range('E73').Value = 'papryka czerwona'
$ws.Range('H73').Value = 'kilka plastrów'
$ws.Range('E74').Value = 'sałata'
$ws.Range('H74').Value = '1 liść'
$ws.Range('E75').Value = 'mięso z piersi z kurczaka'
$ws.Range('E76').Value = 'cebulą '
$ws.Range('H76').Value = '1/4 szt'
$ws.Range('E77').Value = 'jajko'
$ws.Range('H77').Value = '1 szt.'
$ws.Range('E78').Value = 'bułka'
$ws.Range('H78').Value = '1/3 szt.'
$ws.Range('E79').Value = 'X'
$ws.Range('B80').Value = 'obiad'
$ws.Range('D80').Value = 'Roladki schabowe z ziemniakami z pietruszką, Surówka ulubiona'
$ws.Range('H80').Value = '1 plaster'
$ws.Range('I80').Value = 'Schab rozbić na cienki kotlet, doprawić przyprawami; do środka włożyć ser żółty i paprykę ; zawinąć w roladkę; folię aluminiową bardzo delikatnie posmarować olejem; zawinąć w nią roladkę i zapiekać w piekarniku nagrzanym do 220°C, około 30-40 min; podawać z ugotowanymi ziemniakami posypanymi pietruszką i surówką.'
$ws.Range('E81').Value = 'ser żółty'
$ws.Range('H81').Value = '1 plaster'
$ws.Range('E82').Value = 'papryka czerwona'
$ws.Range('H82').Value = '2 plasterki'
$ws.Range('E83').Value = 'papryka słodka'
$ws.Range('H83').Value = '2 plasterki'
$ws.Range('H84').Value = '2 szt.'
$ws.Range('E85').Value = 'pietruszka natka'
$ws.Range('H85').Value = '1 szt.'
$ws.Range('E86').Value = 'surówka gotowa'
$ws.Range('H86').Value = '1/2 opakowania'
$ws.Range('E87').Value = 'sok pomidorowy'
$ws.Range('H87').Value = 'szklanka'
$ws.Range('E88').Value = 'X'
$ws.Range('B89').Value = 'kolacja'
$ws.Range('D89').Value = 'Bułka pełnoziarnistaz twarożkiem ze szczypiorkiem'
$ws.Range('E89').Value = 'chleb żytni razowy'
$ws.Range('H89').Value = '1 kromka'
$ws.Range('I89').Value = 'Wymieszać twarożek ze szczypiorkiem i solą; kanapki przygotować z serkiem '
$ws.Range('E90').Value = 'masło'
$ws.Range('H90').Value = '1/2 lyzeczki'
$ws.Range('I90').Value = 'i szynką'
$ws.Range('E91').Value = 'twaróg chudy'
$ws.Range('E92').Value = 'mleko'
$ws.Range('E93').Value = 'szczypiorek '
$ws.Range('H94').Value = '2 plasterki'
$ws.Range('E95').Value = 'ogórek kiszony/konserwowy'
$ws.Range('H95').Value = '1 szt.'
$ws.Range('E96').Value = 'X'
$ws.Range('B97').Value = 'kolacja'
$ws.Range('D97').Value = 'Kanapki z pasztetem i serem'
$ws.Range('E97').Value = 'chleb żytni razowy'
$ws.Range('H97').Value = '2 kromki'
$ws.Range('I97').Value = 'Przygotować kanapki: jedna: pasztet + papryka + sałata, druga: ser żółty + keczup.'
$ws.Range('E98').Value = 'pasztet z kurczaka'
$ws.Range('H98').Value = 'łyżka'
$ws.Range('E99').Value = 'papryka czerwona'
$ws.Range('H99').Value = '2 plastry'
$ws.Range('H100').Value = '1 cienki plaster'
$ws.Range('E101').Value = 'sałata'
$ws.Range('H101').Value = '1 lisc'
$ws.Range('H102').Value = '1 łyzeczka'
$ws.Range('E103').Value = 'X'
$ws.Range('B104').Value = 'kolacja'
$ws.Range('D104').Value = 'Papryka faszerowana'
$ws.Range('H104').Value = '1/5 opak'
$ws.Range('I104').Value = 'ugotować ryż; wymieszać z mięsem mielonym z indyka, białkiem jaja i natką pietruszki; przyprawić solą i pieprzem; farszem wypełnić paprykę czerwoną; podlać 1 łyżką wody/bulionu; zapiekać w piekarniku nagrzanym do 200 stopni ok. 30-40 min; 10 min przed końcem posypać starym żółtym serem'
$ws.Range('E105').Value = 'mięso mielone z kurczaka lub indyka'
$ws.Range('H105').Value = '1/3 szkl.'
$ws.Range('E106').Value = 'jajko'
$ws.Range('H106').Value = '1 szt.'
$ws.Range('E107').Value = 'papryka czerwona'
$ws.Range('H107').Value = '1 szt.'
$ws.Range('E108').Value = 'zółty ser'
$ws.Range('H108').Value = '1 plasterek'
$ws.Range('E109').Value = 'X'

# --- Sheet view: scroll position reset (drop frozen topLeftCell scroll) and selection moved.
$ws.Range('E108').Select()

# --- New column width for column E (ingredient amount / unit column).
$ws.Columns.Item(5).ColumnWidth = 22.5
